# Updated cryptos list on Sun Feb 11 13:57:47 UTC 2024 with GitHub Actions
# Applies the per-row cell value updates (Price / Volume(1h) columns) and
# the Cosmos/Toncoin row swap (rows 28-29) described by the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'48.077.19"
$ws.Range("E2").Value = "'  +1.58%  "
$ws.Range("D3").Value = "'2.512.30"
$ws.Range("E3").Value = "'  +0.95%  "
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("D5").Value = "'321.49"
$ws.Range("E5").Value = "'  +0.04%  "
$ws.Range("D6").Value = "'109.11"
$ws.Range("E6").Value = "'  +1.14%  "
$ws.Range("E7").Value = "'  +1.61%  "
$ws.Range("E8").Value = "'  +0.10%  "
$ws.Range("D9").Value = "'0.548"
$ws.Range("E9").Value = "'  +2.77%  "
$ws.Range("D10").Value = "'39.86"
$ws.Range("E10").Value = "'  +3.19%  "
$ws.Range("D11").Value = "'20.20"
$ws.Range("E11").Value = "'  +10.01%  "
$ws.Range("E12").Value = "'  +1.09%  "
$ws.Range("E13").Value = "'  +0.96%  "
$ws.Range("D14").Value = "'7.21"
$ws.Range("E14").Value = "'  +1.45%  "
$ws.Range("D15").Value = "'2.910.25"
$ws.Range("E15").Value = "'  +1.10%  "
$ws.Range("D16").Value = "'2.520.82"
$ws.Range("E16").Value = "'  +1.60%  "
$ws.Range("D17").Value = "'0.849"
$ws.Range("E17").Value = "'  +0.44%  "
$ws.Range("D18").Value = "'47.937.41"
$ws.Range("E18").Value = "'  +1.50%  "
$ws.Range("D19").Value = "'13.18"
$ws.Range("E19").Value = "'  +3.02%  "
$ws.Range("D20").Value = "'6.59"
$ws.Range("E20").Value = "'  -0.21%  "
$ws.Range("D21").Value = "'0.0₃0945"
$ws.Range("E21").Value = "'  +1.48%  "
$ws.Range("E22").Value = "'  -0.28%  "
$ws.Range("D23").Value = "'71.87"
$ws.Range("E23").Value = "'  +2.32%  "
$ws.Range("D24").Value = "'276.25"
$ws.Range("E24").Value = "'  +12.56%  "
$ws.Range("D25").Value = "'2.57"
$ws.Range("E25").Value = "'  +0.70%  "
$ws.Range("D27").Value = "'25.91"
$ws.Range("E27").Value = "'  +0.92%  "
$ws.Range("B28").Value = "'Cosmos"
$ws.Range("C28").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'10.11"
$ws.Range("E28").Value = "'  +1.31%  "
$ws.Range("B29").Value = "'Toncoin"
$ws.Range("C29").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.25"
$ws.Range("E29").Value = "'  +3.03%  "
$ws.Range("E30").Value = "'  +5.95%  "
$ws.Range("D31").Value = "'35.41"
$ws.Range("E31").Value = "'  +2.44%  "
$ws.Range("D32").Value = "'49.53"
$ws.Range("E32").Value = "'  -0.04%  "
$ws.Range("D33").Value = "'19.58"
$ws.Range("E33").Value = "'  -5.09%  "
$ws.Range("D34").Value = "'5.36"
$ws.Range("E34").Value = "'  +0.85%  "
$ws.Range("D35").Value = "'1.01"
$ws.Range("E35").Value = "'  -0.01%  "
$ws.Range("D36").Value = "'0.0782"
$ws.Range("E36").Value = "'  +0.20%  "
$ws.Range("E37").Value = "'  +0.60%  "
$ws.Range("D38").Value = "'4.66"
$ws.Range("E38").Value = "'  +0.08%  "
$ws.Range("E39").Value = "'  +1.84%  "
$ws.Range("D40").Value = "'123.12"
$ws.Range("E40").Value = "'  +3.97%  "
$ws.Range("E41").Value = "'  +0.68%  "
$ws.Range("E42").Value = "'  -0.19%  "
$ws.Range("D43").Value = "'21.83"
$ws.Range("E43").Value = "'  -3.62%  "
$ws.Range("D44").Value = "'0.0299"
$ws.Range("E44").Value = "'  +1.34%  "
$ws.Range("D45").Value = "'2.033.05"
$ws.Range("E45").Value = "'  +2.45%  "
$ws.Range("D46").Value = "'3.11"
$ws.Range("E46").Value = "'  +2.55%  "
$ws.Range("D47").Value = "'1.86"
$ws.Range("E47").Value = "'  +5.59%  "
$ws.Range("D48").Value = "'1.99"
$ws.Range("E48").Value = "'  -0.41%  "
$ws.Range("D49").Value = "'9.02"
$ws.Range("E49").Value = "'  -0.38%  "
$ws.Range("D50").Value = "'5.17"
$ws.Range("E50").Value = "'  +1.45%  "
$ws.Range("D51").Value = "'79.91"
$ws.Range("E51").Value = "'  +3.71%  "
